$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = "Índice"
$ws.Cells.Item(1,2).Value = "Distancia"
$ws.Cells.Item(1,3).Value = "max"
$ws.Cells.Item(1,4).Value = "min"
$ws.Cells.Item(1,5).Value = "Tempo"
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 2242.833333333333
$ws.Cells.Item(2,3).Value = 2376
$ws.Cells.Item(2,4).Value = 2098
$ws.Cells.Item(2,5).Value = 0.02966042359670003
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 2247.266666666667
$ws.Cells.Item(3,3).Value = 2331
$ws.Cells.Item(3,4).Value = 2165
$ws.Cells.Item(3,5).Value = 0.0296846866607666
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 1987.2
$ws.Cells.Item(4,3).Value = 2254
$ws.Cells.Item(4,4).Value = 1784
$ws.Cells.Item(4,5).Value = 0.03261879285176595
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 2313.8
$ws.Cells.Item(5,3).Value = 2447
$ws.Cells.Item(5,4).Value = 2199
$ws.Cells.Item(5,5).Value = 0.03328153292338053
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 1974.966666666667
$ws.Cells.Item(6,3).Value = 2133
$ws.Cells.Item(6,4).Value = 1768
$ws.Cells.Item(6,5).Value = 0.03264238039652507
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = 2146.7
$ws.Cells.Item(7,3).Value = 2266
$ws.Cells.Item(7,4).Value = 2051
$ws.Cells.Item(7,5).Value = 0.03248660564422608
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = 2228.233333333333
$ws.Cells.Item(8,3).Value = 2355
$ws.Cells.Item(8,4).Value = 2136
$ws.Cells.Item(8,5).Value = 0.03297599156697591
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = 2261.033333333333
$ws.Cells.Item(9,3).Value = 2476
$ws.Cells.Item(9,4).Value = 2100
$ws.Cells.Item(9,5).Value = 0.03347345987955729
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 2425.133333333333
$ws.Cells.Item(10,3).Value = 2588
$ws.Cells.Item(10,4).Value = 2267
$ws.Cells.Item(10,5).Value = 0.03031268914540609
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 1968.733333333333
$ws.Cells.Item(11,3).Value = 2103
$ws.Cells.Item(11,4).Value = 1938
$ws.Cells.Item(11,5).Value = 0.03040310541788737
